$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: N (sample size), chi-square critical values --------------------
$ws.Rows.Item(6).RowHeight = 17

$ws.Range("K6").Value = "crit upper"
$ws.Range("L6").Formula = "=_xlfn.CHISQ.INV(1 - 0.025, P6-1)"
$ws.Range("L6").Font.Size = 10
$ws.Range("L6").Font.Name = "Arial Unicode MS"

$ws.Range("M6").Value = "crit lower"
$ws.Range("N6").Formula = "=_xlfn.CHISQ.INV( 0.025, P6-1)"
$ws.Range("N6").Font.Size = 10
$ws.Range("N6").Font.Name = "Arial Unicode MS"

$ws.Range("O6").Value = "N"
$ws.Range("P6").Value = 14

# --- Row 7: variance confidence interval (chi-square) -----------------------
$ws.Range("K7").Value = "lower bd"
$ws.Range("L7").Formula = "=(P6-1)*P5/L6"
$ws.Range("M7").Value = "upper bd"
$ws.Range("N7").Formula = "=(P6-1)*P5/N6"

# --- Row 8: mean of log nu ---------------------------------------------------
$ws.Range("O8").Value = "mean (log nu)"
$ws.Range("P8").Formula = "=AVERAGE(A4:N4)"

# --- Row 9: t critical value, std dev of log nu -----------------------------
$ws.Range("K9").Value = "crit"
$ws.Range("L9").Formula = "=_xlfn.T.INV.2T(0.025, P6-1)"
$ws.Range("M9").Value = "std dev"
$ws.Range("N9").Formula = "=SQRT(P5)"

# --- Row 10: mean (log nu) confidence interval ------------------------------
$ws.Range("K10").Value = "lower bd"
$ws.Range("L10").Formula = "=P8-(L9*N9/SQRT(P6))"
$ws.Range("M10").Value = "upper bd"
$ws.Range("N10").Formula = "=P8+(L9*N9/SQRT(P6))"

# --- Row 11: mean of nu ------------------------------------------------------
$ws.Range("O11").Value = "mean (nu)"
$ws.Range("P11").Formula = "=AVERAGE(G3:N3)"

# --- Row 12: std dev / variance of nu ----------------------------------------
$ws.Range("M12").Value = "std dev"
$ws.Range("N12").Formula = "=SQRT(P12)"
$ws.Range("O12").Value = "variance"
$ws.Range("P12").Formula = "=VARA(A3:N3)"

# --- Row 13: mean (nu) confidence interval -----------------------------------
$ws.Range("K13").Value = "lower bd"
$ws.Range("L13").Formula = "=P11-(L9*N12/SQRT(P6))"
$ws.Range("M13").Value = "upper bd"
$ws.Range("N13").Formula = "=P11+(L9*N12/SQRT(P6))"

# --- Row 14: nu^5 for each observation + average -----------------------------
$ws.Range("A14").Formula = "=A3^5"
$ws.Range("B14").Formula = "=B3^5"
$ws.Range("C14").Formula = "=C3^5"
$ws.Range("D14").Formula = "=D3^5"
$ws.Range("E14").Formula = "=E3^5"
$ws.Range("F14").Formula = "=F3^5"
$ws.Range("G14").Formula = "=G3^5"
$ws.Range("H14").Formula = "=H3^5"
$ws.Range("I14").Formula = "=I3^5"
$ws.Range("J14").Formula = "=J3^5"
$ws.Range("K14").Formula = "=K3^5"
$ws.Range("L14").Formula = "=L3^5"
$ws.Range("M14").Formula = "=M3^5"
$ws.Range("N14").Formula = "=N3^5"
$ws.Range("O14").Value = "nu^5"
$ws.Range("P14").Formula = "=AVERAGE(A14:N14)"

# --- Row 15: std dev / variance of nu^5 --------------------------------------
$ws.Range("M15").Value = "std dev"
$ws.Range("N15").Formula = "=SQRT(P15)"
$ws.Range("O15").Value = "variance"
$ws.Range("P15").Formula = "=VARA(A14:N14)"

# --- Row 16: mean (nu^5) confidence interval ---------------------------------
$ws.Range("K16").Value = "lower bd"
$ws.Range("L16").Formula = "=P14-(L9*N15/SQRT(P6))"
$ws.Range("M16").Value = "upper bd"
$ws.Range("N16").Formula = "=P14+(L9*N15/SQRT(P6))"

# --- View state: scroll / selection matches the saved file ------------------
$ws.Range("L16").Select()
